$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.234.86'
$ws.Range('E2').Value = '  +0.25%  '

$ws.Range('D3').Value = '1.788.06'
$ws.Range('E3').Value = '  -0.29%  '

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').Value = "'225.78"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.77%  '

$ws.Range('D6').Value = "'0.555"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.43%  '

$ws.Range('E7').Value = '  +0.16%  '

$ws.Range('D8').Value = "'32.24"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.21%  '

$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('D10').Value = "'0.0688"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.85%  '

$ws.Range('D11').Value = "'0.0948"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.73%  '

$ws.Range('D12').Value = '2.045.15'

$ws.Range('D13').Value = '1.787.48'
$ws.Range('E13').Value = '  -0.37%  '

$ws.Range('D14').Value = "'10.98"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -5.19%  '

$ws.Range('E15').Value = '  +0.38%  '

$ws.Range('D16').Value = '34.194.84'
$ws.Range('E16').Value = '  +0.19%  '

$ws.Range('D17').Value = "'4.19"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.32%  '

$ws.Range('D18').Value = "'67.95"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.11%  '

$ws.Range('D19').Value = '0.0₃0803'
$ws.Range('E19').Value = '  +2.39%  '

$ws.Range('D20').Value = "'246.31"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.37%  '

$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').Value = "'10.94"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.01%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = "'1.00"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.19%  '

$ws.Range('D23').Value = "'4.17"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.15%  '

$ws.Range('E24').Value = '  +0.22%  '

$ws.Range('D25').Value = "'162.17"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.20%  '

$ws.Range('D26').Value = "'7.17"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.16%  '

$ws.Range('D27').Value = "'16.35"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.20%  '

$ws.Range('E28').Value = '  +1.26%  '

$ws.Range('E29').Value = '  +0.30%  '

$ws.Range('E30').Value = '  -0.80%  '

$ws.Range('E31').Value = '  -0.10%  '

$ws.Range('E32').Value = '  +1.95%  '

$ws.Range('D33').Value = "'3.79"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.40%  '

$ws.Range('E34').Value = '  -1.58%  '

$ws.Range('D35').Value = '1.439.64'
$ws.Range('E35').Value = '  -0.54%  '

$ws.Range('D36').Value = "'2.62"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +9.96%  '

$ws.Range('E37').Value = '  +1.85%  '

$ws.Range('E38').Value = '  +1.44%  '

$ws.Range('D39').Value = "'0.0190"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.32%  '

$ws.Range('D40').Value = "'81.84"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.18%  '

$ws.Range('E41').Value = '  +2.06%  '

$ws.Range('D42').Value = "'14.06"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.40%  '

$ws.Range('E43').Value = '  +1.35%  '

$ws.Range('D44').Value = "'0.922"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.76%  '

$ws.Range('E45').Value = '  +2.15%  '

$ws.Range('D46').Value = "'6.06"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.25%  '

$ws.Range('E47').Value = '  +0.60%  '

$ws.Range('D48').Value = '1.940.66'
$ws.Range('E48').Value = '  -0.53%  '

$ws.Range('D49').Value = "'105.56"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.30%  '

$ws.Range('E50').Value = '  +0.19%  '

$ws.Range('E51').Value = '  -6.45%  '
